$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Move the existing STATUS column (column B) into column C,
#    preserving its original formatting (font/style).
# ---------------------------------------------------------------
$ws.Range("B1:B4").Copy()
$ws.Range("C1:C4").PasteSpecial(-4104)   # xlPasteAll

# Row 4 status text changes from "Doing" to "Done" as part of this edit
$ws.Range("C4").Value = "Done"

# ---------------------------------------------------------------
# 2. Column B becomes a new "Observation" column.
#    B1 header takes on the same look as A1.
# ---------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("B1").Value = "Observation"

# ---------------------------------------------------------------
# 3. New observation rows (5-7) content
# ---------------------------------------------------------------
$ws.Range("A5").Value = "COW table"
$ws.Range("B5").Value = "Update a record will create a new data file and the original file will be kept as it is. Original file size : 1 KB and After modifying file size : 2KB"
$ws.Range("C5").Value = "Done"
$ws.Range("B6").Value = "Deleting records where class_year < 2020. It deleted all the records from year 2019. No records from 2019 are shown in the table, but the data file containing all records from 2019 is present. It will be used to see the historic data."
$ws.Range("B7").Value = "Upserting data from new_table to students"

# ---------------------------------------------------------------
# 4. Formatting for rows 5-7.
#    Build each combined alignment format once on a scratch cell
#    (which already carries the row's base font) and paste the
#    resulting format onto the real cells - this avoids leaving
#    unused intermediate styles behind in the stylesheet.
# ---------------------------------------------------------------
$ws.Range("B5").WrapText = $true
$ws.Range("B6").WrapText = $true

$ws.Range("Z2").VerticalAlignment = -4108     # xlCenter
$ws.Range("Z2").HorizontalAlignment = -4108   # xlCenter
$ws.Range("Z2").Copy()
$ws.Range("A5:A7").PasteSpecial(-4122)        # xlPasteFormats
$ws.Range("C5:C7").PasteSpecial(-4122)        # xlPasteFormats
$ws.Range("Z2").Clear()

$ws.Range("Z2").VerticalAlignment = -4108     # xlCenter
$ws.Range("Z2").HorizontalAlignment = -4131   # xlLeft
$ws.Range("Z2").Copy()
$ws.Range("B7").PasteSpecial(-4122)           # xlPasteFormats
$ws.Range("Z2").Clear()

$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 100

# ---------------------------------------------------------------
# 5. Merge the A and C observation cells across rows 5-7
# ---------------------------------------------------------------
$ws.Range("A5:A7").Merge()
$ws.Range("C5:C7").Merge()

# ---------------------------------------------------------------
# 6. Column widths: A & B share the original column-A width,
#    C gets the original column-B (STATUS) width.
# ---------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 46.33
$ws.Columns.Item(3).ColumnWidth = 31.65

# ---------------------------------------------------------------
# 7. Selection shown when the file was saved
# ---------------------------------------------------------------
[void]$ws.Range("C12").Select()
